# V1.3 Color Display in Excel
# Adds green highlight fill for all TRUE boolean cells, fixes row 49,
# and extends the data pattern through row 73.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix row 49 to follow the repeating pattern (C,E,F,H => True) ---
$ws.Range("C49").Value = $true
$ws.Range("E49").Value = $true
$ws.Range("F49").Value = $true
$ws.Range("H49").Value = $true

# --- Append new data rows 50-73, continuing the same pattern ---
$newRows = @(
    @(50, 49, $false, $false, $false, $false, $false, $false, $false, $false),
    @(51, 50, $false, $false, $false, $false, $false, $false, $false, $false),
    @(52, 51, $false, $true, $false, $true, $true, $false, $true, $false),
    @(53, 52, $false, $false, $false, $false, $false, $false, $false, $false),
    @(54, 53, $false, $false, $false, $false, $false, $false, $false, $false),
    @(55, 54, $false, $true, $true, $true, $true, $true, $true, $false),
    @(56, 55, $false, $false, $false, $false, $false, $false, $false, $false),
    @(57, 56, $false, $false, $false, $false, $false, $false, $false, $false),
    @(58, 57, $false, $true, $false, $true, $true, $false, $true, $false),
    @(59, 58, $false, $false, $false, $false, $false, $false, $false, $false),
    @(60, 59, $false, $false, $false, $false, $false, $false, $false, $false),
    @(61, 60, $false, $true, $false, $true, $true, $false, $true, $false),
    @(62, 61, $false, $false, $false, $false, $false, $false, $false, $false),
    @(63, 62, $false, $false, $false, $false, $false, $false, $false, $false),
    @(64, 63, $false, $true, $true, $true, $true, $true, $true, $false),
    @(65, 64, $false, $false, $false, $false, $false, $false, $false, $false),
    @(66, 65, $false, $false, $false, $false, $false, $false, $false, $false),
    @(67, 66, $false, $true, $false, $true, $true, $false, $true, $false),
    @(68, 67, $false, $false, $false, $false, $false, $false, $false, $false),
    @(69, 68, $false, $false, $false, $false, $false, $false, $false, $false),
    @(70, 69, $false, $true, $false, $true, $true, $false, $true, $false),
    @(71, 70, $false, $false, $false, $false, $false, $false, $false, $false),
    @(72, 71, $false, $false, $false, $false, $false, $false, $false, $false),
    @(73, 72, $false, $false, $false, $false, $false, $false, $false, $false)
)

foreach ($rowDef in $newRows) {
    $r = $rowDef[0]
    for ($col = 1; $col -le 9; $col++) {
        $ws.Cells.Item($r, $col).Value = $rowDef[$col]
    }
}

# --- Apply solid green fill to every cell whose boolean value is TRUE ---
$greenCells = @(
    "C4", "C7", "E7", "H7", "C10", "D10", "E10", "F10", "G10", "H10",
    "C13", "E13", "F13", "H13", "C16", "E16", "F16", "H16", "C19", "D19",
    "E19", "F19", "G19", "H19", "C22", "E22", "F22", "H22", "C25", "E25",
    "F25", "H25", "C28", "D28", "E28", "F28", "G28", "H28", "C31", "E31",
    "F31", "H31", "C34", "E34", "F34", "H34", "C37", "D37", "E37", "F37",
    "G37", "H37", "C40", "E40", "F40", "H40", "C43", "E43", "F43", "H43",
    "C46", "D46", "E46", "F46", "G46", "H46", "C49", "E49", "F49", "H49",
    "C52", "E52", "F52", "H52", "C55", "D55", "E55", "F55", "G55", "H55",
    "C58", "E58", "F58", "H58", "C61", "E61", "F61", "H61", "C64", "D64",
    "E64", "F64", "G64", "H64", "C67", "E67", "F67", "H67", "C70", "E70",
    "F70", "H70"
)

foreach ($ref in $greenCells) {
    $ws.Range($ref).Interior.Color = 65280
}

Write-Output "Applied green fill to $($greenCells.Count) cells; added $($newRows.Count) new rows."
